$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = -21.943
$ws.Range("C11").Value = -12.7982
$ws.Range("A12").Value = -21.36069999999998
$ws.Range("A15").Value = -21.664
$ws.Range("C23").Value = -12.39749999999999
$ws.Range("A27").Value = -21.8254
$ws.Range("A28").Value = -21.9517
$ws.Range("C28").Value = -13.5608
$ws.Range("A31").Value = -21.50949999999999
$ws.Range("A32").Value = -21.8244
$ws.Range("C32").Value = -12.6559
$ws.Range("C34").Value = -11.64180000000001
$ws.Range("A36").Value = -19.30819999999999
$ws.Range("C36").Value = -12.66430000000001
$ws.Range("C37").Value = -12.7182
$ws.Range("A38").Value = -19.47799999999998
$ws.Range("C42").Value = -12.6402
$ws.Range("A46").Value = -21.8288
$ws.Range("C49").Value = -14.10409999999999
$ws.Range("A54").Value = -21.96939999999998
$ws.Range("C54").Value = -12.6255
$ws.Range("A55").Value = -22.02960000000001
$ws.Range("A56").Value = -22.24010000000001
$ws.Range("A67").Value = -21.46139999999998
$ws.Range("A69").Value = -21.65099999999997
$ws.Range("A72").Value = -21.6756
$ws.Range("A73").Value = -19.9934
$ws.Range("C78").Value = -13.5108
$ws.Range("C80").Value = -12.5452
$ws.Range("A83").Value = -21.72049999999999
$ws.Range("A86").Value = -22.12240000000001
$ws.Range("A91").Value = -20.25719999999998
$ws.Range("A93").Value = -21.35510000000001
$ws.Range("C97").Value = -11.54
$ws.Range("A99").Value = -21.85749999999999
$ws.Range("C99").Value = -13.18319999999999
$ws.Range("C100").Value = -11.9002
$ws.Range("C101").Value = -12.9982
$ws.Range("A104").Value = -21.28389999999999
$ws.Range("A105").Value = -19.62219999999998
